$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 to 6 (keep header row 1 and first data row, which will become the RUG514 row)
$ws.Range("A3:E6").EntireRow.Delete()

# Delete column C (old "max" column) - shifts old D,E into C,D
$ws.Range("C1:C2").EntireColumn.Delete()

# Update row 2 data to reflect RUG514.fasta values
$ws.Range("A2").Value = "RUG514.fasta"
$ws.Range("B2").Value = 46250.81325848302
$ws.Range("C2").Value = "o__Desulfovibrionales"
$ws.Range("D2").Value = "o__Desulfovibrionales"
